# Apply data updates to Sheet1: repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "dSF" column (F) values for the affected rows
$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -6
$ws.Range("F10").Value = -3
$ws.Range("F12").Value = -2
